# edit.ps1 - apply "Alternative file renamed as gandhi" change
#
# 1. Append, to the first paragraph ("This is a Microsoft word
#    document."), a two-space run followed by a dark-red
#    "(This is a change - Version for branch alternate)" annotation
#    (uses an en dash, split across three runs, matching the source edit).
# 2. Add a new, empty trailing paragraph just before the final
#    section break.

$d = $word.ActiveDocument

# --- 1. Annotate the first paragraph -----------------------------------

$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
# Exclude the paragraph mark itself so new runs land inside paragraph 1,
# not after it.
$r.End = $r.End - 1
$r.Collapse(0)

# Plain two-space run (inherits the surrounding - unformatted - style).
$r.InsertAfter("  ")
$r.Collapse(0)

# Dark red (C00000) annotation, split into three runs as in the source.
$r.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r.Font.Color = 192
$r.Collapse(0)

$r.InsertAfter("rsion for branch alternate")
$r.Font.Color = 192
$r.Collapse(0)

$r.InsertAfter(")")
$r.Font.Color = 192
$r.Collapse(0)

# --- 2. Trailing empty paragraph ----------------------------------------

$endR = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endR.InsertParagraphAfter()
$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Range.Style = "Normal"
